$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 keeps its original "Sending/Ligand/Receptor/Target cluster" labels
# (D2 stays "MuSCs"); only the numeric NATMI statistics are refreshed with
# the new TPM-derived values, splitting the previous combined
# "Receptor-expressing cells" count (2) down to the MuSCs-only figures.
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.003092666666666667
$ws.Range("N2").Value = 0.009278
$ws.Range("O2").Value = 0.03934390080485799
$ws.Range("P2").Value = 0.03934390080485798
$ws.Range("Q2").Value = 0.0002204061062222222
$ws.Range("R2").Value = 0.001983654956
$ws.Range("S2").Value = 0.03934390080485799
$ws.Range("T2").Value = 0.03934390080485798

# New row 3 - same sending cluster / ligand / receptor, new target cluster "ECs"
$ws.Range("A3").Value = "FAPs"
$ws.Range("B3").Value = "Wnt1"
$ws.Range("C3").Value = "Fzd10"
$ws.Range("D3").Value = "ECs"
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.07126733333333334
$ws.Range("H3").Value = 0.213802
$ws.Range("I3").Value = 1
$ws.Range("J3").Value = 1
$ws.Range("K3").Value = 2
$ws.Range("L3").Value = 0.6666666666666666
$ws.Range("M3").Value = 0.07551333333333334
$ws.Range("N3").Value = 0.22654
$ws.Range("O3").Value = 0.9606560991951421
$ws.Range("P3").Value = 0.9606560991951419
$ws.Range("Q3").Value = 0.005381633897777778
$ws.Range("R3").Value = 0.04843470507999999
$ws.Range("S3").Value = 0.9606560991951421
$ws.Range("T3").Value = 0.9606560991951419

$wb.Save()
